$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 133.41667
$ws.Range("I55").Value = 101.833336
$ws.Range("J55").Value = 165
$ws.Range("K55").Value = 101.833336
$ws.Range("L55").Value = 165
$ws.Range("M55").Value = 112.166664
$ws.Range("N55").Value = -593

$ws.Range("H100").Value = 1678.6666
$ws.Range("I100").Value = 1410
$ws.Range("J100").Value = 2100.8572
$ws.Range("K100").Value = 1410
$ws.Range("L100").Value = 2100.8572
$ws.Range("M100").Value = -869
$ws.Range("N100").Value = -3182.8572

$ws.Range("H112").Value = 1905.9259
$ws.Range("J112").Value = 2098.261
$ws.Range("L112").Value = 6294.782999999999
$ws.Range("N112").Value = -8510.782999999999

$ws.Range("H135").Value = 19234.123
$ws.Range("I135").Value = 23784.887
$ws.Range("J135").Value = 3831.5386
$ws.Range("K135").Value = 214063.983
$ws.Range("L135").Value = 34483.8474
$ws.Range("M135").Value = -211528.983
$ws.Range("N135").Value = -39553.8474

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2486.42
$ws.Range("I32").Value = 1955.3507
$ws.Range("K32").Value = 1955.3507
$ws.Range("M32").Value = -1668.3507

$ws.Range("H61").Value = 17893818
$ws.Range("I61").Value = 18538024
$ws.Range("J61").Value = 500250
$ws.Range("K61").Value = 18538024
$ws.Range("L61").Value = 500250
$ws.Range("M61").Value = -18537812
$ws.Range("N61").Value = -500674

$ws.Range("H74").Value = 6707839
$ws.Range("I74").Value = 8573611
$ws.Range("J74").Value = 92829.55
$ws.Range("K74").Value = 8573611
$ws.Range("L74").Value = 92829.55
$ws.Range("M74").Value = -8572737
$ws.Range("N74").Value = -94577.55

$ws.Range("H77").Value = 6707839
$ws.Range("I77").Value = 8573611
$ws.Range("J77").Value = 92829.55
$ws.Range("K77").Value = 42868055
$ws.Range("L77").Value = 464147.75
$ws.Range("M77").Value = -42863687
$ws.Range("N77").Value = -472883.75

$ws.Range("H136").Value = 17893818
$ws.Range("I136").Value = 18538024
$ws.Range("J136").Value = 500250
$ws.Range("K136").Value = 55614072
$ws.Range("L136").Value = 1500750
$ws.Range("M136").Value = -55611522
$ws.Range("N136").Value = -1505850

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1032.3636
$ws.Range("I99").Value = 916.6667
$ws.Range("J99").Value = 1171.2
$ws.Range("K99").Value = 916.6667
$ws.Range("L99").Value = 1171.2
$ws.Range("M99").Value = 581.3333
$ws.Range("N99").Value = -4167.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 1340.8
$ws.Range("I10").Value = 234.66667
$ws.Range("J10").Value = 3000
$ws.Range("K10").Value = 234.66667
$ws.Range("L10").Value = 3000
$ws.Range("M10").Value = -95.66667000000001
$ws.Range("N10").Value = -3278

$ws.Range("H16").Value = 83335340
$ws.Range("J16").Value = 111113120
$ws.Range("L16").Value = 111113120
$ws.Range("N16").Value = -111113694

$ws.Range("H58").Value = 18183138
$ws.Range("I58").Value = 20409344
$ws.Range("J58").Value = 2450.5
$ws.Range("K58").Value = 20409344
$ws.Range("L58").Value = 2450.5
$ws.Range("M58").Value = -20409141
$ws.Range("N58").Value = -2856.5

$ws.Range("H74").Value = 25301.44
$ws.Range("J74").Value = 25301.44
$ws.Range("L74").Value = 25301.44
$ws.Range("N74").Value = -27049.44

$ws.Range("H77").Value = 25301.44
$ws.Range("J77").Value = 25301.44
$ws.Range("L77").Value = 75904.31999999999
$ws.Range("N77").Value = -84640.31999999999

$ws.Range("H113").Value = 83335340
$ws.Range("J113").Value = 111113120
$ws.Range("L113").Value = 111113120
$ws.Range("N113").Value = -111117460

$ws.Range("H122").Value = 2321.125
$ws.Range("I122").Value = 1922.4
$ws.Range("J122").Value = 2985.6667
$ws.Range("K122").Value = 5767.200000000001
$ws.Range("L122").Value = 8957.000100000001
$ws.Range("M122").Value = -3317.200000000001
$ws.Range("N122").Value = -13857.0001

$ws.Range("H136").Value = 18183138
$ws.Range("I136").Value = 20409344
$ws.Range("J136").Value = 2450.5
$ws.Range("K136").Value = 61228032
$ws.Range("L136").Value = 7351.5
$ws.Range("M136").Value = -61225482
$ws.Range("N136").Value = -12451.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 527.2353000000001
$ws.Range("I68").Value = 410.44446
$ws.Range("J68").Value = 658.625
$ws.Range("K68").Value = 1231.33338
$ws.Range("L68").Value = 1975.875
$ws.Range("M68").Value = -420.33338
$ws.Range("N68").Value = -3597.875

$ws.Range("H71").Value = 527.2353000000001
$ws.Range("I71").Value = 410.44446
$ws.Range("J71").Value = 658.625
$ws.Range("K71").Value = 3694.00014
$ws.Range("L71").Value = 5927.625
$ws.Range("M71").Value = 361.9998599999999
$ws.Range("N71").Value = -14039.625

$ws.Range("H126").Value = 2877.7778
$ws.Range("I126").Value = 1466.6666
$ws.Range("K126").Value = 4399.9998
$ws.Range("M126").Value = 540.0002000000004

$ws.Range("H129").Value = 2690011.8
$ws.Range("I129").Value = 1467.5
$ws.Range("J129").Value = 4388039.5
$ws.Range("K129").Value = 4402.5
$ws.Range("L129").Value = 13164118.5
$ws.Range("M129").Value = 597.5
$ws.Range("N129").Value = -13174118.5

$ws.Range("H131").Value = 1023.55
$ws.Range("J131").Value = 1120.8269
$ws.Range("L131").Value = 3362.4807
$ws.Range("N131").Value = -13442.4807

$ws.Range("H132").Value = 2633.8096
$ws.Range("I132").Value = 1907.1428
$ws.Range("J132").Value = 2997.1428
$ws.Range("K132").Value = 17164.2852
$ws.Range("L132").Value = 26974.2852
$ws.Range("M132").Value = -14634.2852
$ws.Range("N132").Value = -32034.2852

$ws.Range("H136").Value = 2822
$ws.Range("I136").Value = 2524
$ws.Range("J136").Value = 3120
$ws.Range("K136").Value = 7572
$ws.Range("L136").Value = 9360
$ws.Range("M136").Value = -2472
$ws.Range("N136").Value = -19560

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 4648.3335
$ws.Range("I43").Value = 633.3333
$ws.Range("K43").Value = 633.3333
$ws.Range("M43").Value = -482.3333

$ws.Range("H132").Value = 65692.28999999999
$ws.Range("I132").Value = 35567.656
$ws.Range("J132").Value = 502499.5
$ws.Range("K132").Value = 106702.968
$ws.Range("L132").Value = 1507498.5
$ws.Range("M132").Value = -104172.968
$ws.Range("N132").Value = -1512558.5

$ws.Range("H133").Value = 77070
$ws.Range("J133").Value = 95250.5
$ws.Range("L133").Value = 95250.5
$ws.Range("N133").Value = -105370.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 15000
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 15000
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 15000
$ws.Range("N7").Value = -15224
$ws.Range("M7").ClearContents()

$ws.Range("H61").Value = 1162.4166
$ws.Range("I61").Value = 1337
$ws.Range("J61").Value = 918
$ws.Range("K61").Value = 1337
$ws.Range("L61").Value = 918
$ws.Range("M61").Value = -1135
$ws.Range("N61").Value = -1322

$ws.Range("H113").Value = 1162.4166
$ws.Range("I113").Value = 1337
$ws.Range("J113").Value = 918
$ws.Range("K113").Value = 1337
$ws.Range("L113").Value = 918
$ws.Range("M113").Value = 833
$ws.Range("N113").Value = -5258

$ws.Range("H126").Value = 15000
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 15000
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 45000
$ws.Range("N126").Value = -49940
$ws.Range("M126").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 34224.098
$ws.Range("I81").Value = 1487.8462
$ws.Range("J81").Value = 57866.945
$ws.Range("K81").Value = 2975.6924
$ws.Range("L81").Value = 115733.89
$ws.Range("M81").Value = -1914.6924
$ws.Range("N81").Value = -117855.89

$ws.Range("H84").Value = 34224.098
$ws.Range("I84").Value = 1487.8462
$ws.Range("J84").Value = 57866.945
$ws.Range("K84").Value = 14878.462
$ws.Range("L84").Value = 578669.45
$ws.Range("M84").Value = -9574.462
$ws.Range("N84").Value = -589277.45

$ws.Range("H107").Value = 660
$ws.Range("I107").Value = 650
$ws.Range("J107").Value = 666.6667
$ws.Range("K107").Value = 1950
$ws.Range("L107").Value = 2000.0001
$ws.Range("M107").Value = -30
$ws.Range("N107").Value = -5840.0001

$ws.Range("H122").Value = 1810.9487
$ws.Range("I122").Value = 1312.409
$ws.Range("J122").Value = 2456.1177
$ws.Range("K122").Value = 3937.227
$ws.Range("L122").Value = 7368.353099999999
$ws.Range("M122").Value = -1487.227
$ws.Range("N122").Value = -12268.3531

$ws.Range("H136").Value = 27887.855
$ws.Range("I136").Value = 22062.936
$ws.Range("J136").Value = 37328.242
$ws.Range("K136").Value = 66188.808
$ws.Range("L136").Value = 111984.726
$ws.Range("M136").Value = -63638.808
$ws.Range("N136").Value = -117084.726
